{"js": "// Charge grid: the \"Plea\" row becomes \"No Contest\" for all three charges,\n// and the previously-blank \"Finding\" cell for the third charge is filled\n// in with \"Guilty\".\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet pleaRow = null;\nlet findingRow = null;\nfor (const row of table.rows.items) {\n  const label = row.cells.items[0].body.text.trim();\n  if (label === \"Plea\") {\n    pleaRow = row;\n  } else if (label === \"Finding\") {\n    findingRow = row;\n  }\n}\n\n// Load the paragraph for every cell we are about to touch so we can get a\n// precise Range and replace its text in place (preserving run/paragraph\n// formatting instead of rewriting the whole cell body).\nfor (let c = 1; c < pleaRow.cells.items.length; c++) {\n  pleaRow.cells.items[c].body.paragraphs.load(\"items\");\n}\nconst findingLastCell =\n  findingRow.cells.items[findingRow.cells.items.length - 1];\nfindingLastCell.body.paragraphs.load(\"items\");\nawait context.sync();\n\n// Plea columns for all three charges -> \"No Contest\".\nfor (let c = 1; c < pleaRow.cells.items.length; c++) {\n  const para = pleaRow.cells.items[c].body.paragraphs.items[0];\n  para.getRange().insertText(\"No Contest\", \"Replace\");\n}\n\n// Finding column for the third charge (previously blank) -> \"Guilty\".\nconst findingPara = findingLastCell.body.paragraphs.items[0];\nfindingPara.getRange().insertText(\"Guilty\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Charge grid: the \"Plea\" row becomes \"No Contest\" for all three charges,\n# and the previously-blank \"Finding\" cell for the third charge is filled\n# in with \"Guilty\".\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$pleaRow = $null\n$findingRow = $null\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $labelCell = $table.Cell($r, 1)\n    $label = $labelCell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($label -eq \"Plea\") {\n        $pleaRow = $r\n    } elseif ($label -eq \"Finding\") {\n        $findingRow = $r\n    }\n}\n\n# Plea columns for all three charges -> \"No Contest\".\nfor ($c = 2; $c -le $table.Columns.Count; $c++) {\n    $cell = $table.Cell($pleaRow, $c)\n    $cell.Range.Text = \"No Contest\"\n}\n\n# Finding column for the third charge (previously blank) -> \"Guilty\".\n$lastCol = $table.Columns.Count\n$findingCell = $table.Cell($findingRow, $lastCol)\n$findingCell.Range.Text = \"Guilty\"\n"}
